$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.459.19"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.871.63"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'331.10"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4597"
$ws.Range("E7").Value = "  -2.23%  "
$ws.Range("D8").Value = "'0.4044"
$ws.Range("E8").Value = "  +3.03%  "
$ws.Range("D9").Value = "'47.56"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "'0.07859"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "'0.9873"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "'21.41"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "1.866.36"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "'5.845"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "'7.016"
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("D16").Value = "'0.9997"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "'88.36"
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").Value = "'0.06547"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "'0.00001019"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("D20").Value = "'17.21"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").Value = "'0.9993"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "28.448.77"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "'5.342"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "'10.87"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'2.250"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "2.085.67"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "'157.37"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "'19.29"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").Value = "'2.069"
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("D30").Value = "'5.312"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").Value = "'117.37"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "'0.9575"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").Value = "'0.09333"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("D34").Value = "'3.596"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "'1.393"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").Value = "'5.229"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("D37").Value = "'0.06027"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'0.02208"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").Value = "'8.295"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D41").Value = "'0.9995"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "'0.5772"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").Value = "'0.1812"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D45").Value = "'1.242"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("D46").Value = "'2.305"
$ws.Range("E46").Value = "  +15.52%  "
$ws.Range("D48").Value = "'11.82"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "'0.07169"
$ws.Range("E49").Value = "  +4.06%  "
$ws.Range("D50").Value = "'1.883"
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("D51").Value = "'109.87"
$ws.Range("E51").Value = "  -0.95%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E44").Value = "  -3.35%  "
$ws.Range("E47").Value = "  -3.27%  "
